$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '36.510.77'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.31%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.954.46'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.66%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.90'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.51%  '
$ws.Range('E6').Value = '  +0.35%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '58.13'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.58%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.375'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.44%  '
$ws.Range('E10').Value = '  -7.43%  '
$ws.Range('E11').Value = '  -0.53%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.03'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.01%  '
$ws.Range('E13').Value = '  +2.40%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.241.38'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.66%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.17'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.11%  '
$ws.Range('E16').Value = '  +1.82%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.959.23'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.22%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '36.467.83'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.32%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.70'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.49%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0845'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.20%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '229.13'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.03'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.52%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('E24').Value = '  +1.89%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.34'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.78%  '
$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.139'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +6.54%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.12'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.58%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '160.17'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.48%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.32'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.22%  '
$ws.Range('E30').Value = '  +1.55%  '
$ws.Range('E31').Value = '  +4.30%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.71'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.98%  '
$ws.Range('E33').Value = '  -4.17%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.39'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.96%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.46'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +13.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.06%  '
$ws.Range('E37').Value = '  +5.79%  '
$ws.Range('E38').Value = '  -1.30%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.29'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -13.65%  '
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('E41').Value = '  +1.51%  '
$ws.Range('E42').Value = '  +0.06%  '
$ws.Range('E43').Value = '  -0.23%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.371.43'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.24%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '15.73'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '87.77'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.65%  '
$ws.Range('E47').Value = '  -0.30%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.12'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.30%  '
$ws.Range('E49').Value = '  +0.27%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.131.66'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.60%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '43.79'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.54%  '
